$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Remove the row for NM "11.245.954" (Extintor incên. c/carga de pó- 50 Kg (80:BC) *).
# This is row 34 in the sheet (row 1 is the header).
$ws.Rows.Item(34).Delete()
